$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 6716.6665
$ws.Range("I29").Value = 5075
$ws.Range("K29").Value = 15225
$ws.Range("M29").Value = -14944
$ws.Range("H40").Value = 4700
$ws.Range("I40").Value = 4100
$ws.Range("K40").Value = 4100
$ws.Range("M40").Value = -3925
$ws.Range("H86").Value = 5455.625
$ws.Range("J86").Value = 9361.75
$ws.Range("L86").Value = 9361.75
$ws.Range("N86").Value = -11607.75
$ws.Range("H89").Value = 5455.625
$ws.Range("J89").Value = 9361.75
$ws.Range("L89").Value = 46808.75
$ws.Range("N89").Value = -58040.75
$ws.Range("H92").Value = 333.6875
$ws.Range("I92").Value = 407.81818
$ws.Range("K92").Value = 407.81818
$ws.Range("M92").Value = 840.18182
$ws.Range("H98").Value = 2468
$ws.Range("I98").Value = 2519.3333
$ws.Range("J98").Value = 2006
$ws.Range("K98").Value = 2519.3333
$ws.Range("L98").Value = 2006
$ws.Range("M98").Value = -1021.3333
$ws.Range("N98").Value = -5002
$ws.Range("H99").Value = 3223.4285
$ws.Range("I99").Value = 1947.5
$ws.Range("K99").Value = 5842.5
$ws.Range("M99").Value = -4344.5
$ws.Range("H101").Value = 2045
$ws.Range("I101").Value = 938
$ws.Range("J101").Value = 3594.8
$ws.Range("K101").Value = 2814
$ws.Range("L101").Value = 10784.4
$ws.Range("M101").Value = -1192
$ws.Range("N101").Value = -14028.4
$ws.Range("H106").Value = 974.4666999999999
$ws.Range("I106").Value = 972.7143
$ws.Range("J106").Value = 999
$ws.Range("K106").Value = 972.7143
$ws.Range("L106").Value = 999
$ws.Range("M106").Value = -341.7143
$ws.Range("N106").Value = -2261
$ws.Range("H112").Value = 1685.2941
$ws.Range("I112").Value = 1961.5
$ws.Range("J112").Value = 1600.3077
$ws.Range("K112").Value = 5884.5
$ws.Range("L112").Value = 4800.9231
$ws.Range("M112").Value = -4776.5
$ws.Range("N112").Value = -7016.9231
$ws.Range("H122").Value = 2468
$ws.Range("I122").Value = 2519.3333
$ws.Range("J122").Value = 2006
$ws.Range("K122").Value = 7557.999899999999
$ws.Range("L122").Value = 6018
$ws.Range("M122").Value = -5107.999899999999
$ws.Range("N122").Value = -10918
$ws.Range("H132").Value = 64327.562
$ws.Range("I132").Value = 85168.25
$ws.Range("J132").Value = 1805.5
$ws.Range("K132").Value = 255504.75
$ws.Range("L132").Value = 5416.5
$ws.Range("M132").Value = -252974.75
$ws.Range("N132").Value = -10476.5
$ws.Range("H138").Value = 3004.2444
$ws.Range("J138").Value = 3389.8235
$ws.Range("L138").Value = 10169.4705
$ws.Range("N138").Value = -20449.4705
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2562.6843
$ws.Range("I32").Value = 2701.3235
$ws.Range("K32").Value = 2701.3235
$ws.Range("M32").Value = -2414.3235
$ws.Range("H45").Value = 2142.2856
$ws.Range("I45").Value = 2099.2
$ws.Range("K45").Value = 2099.2
$ws.Range("M45").Value = -1722.2
$ws.Range("H61").Value = 4866.7856
$ws.Range("I61").Value = 6483
$ws.Range("J61").Value = 4426
$ws.Range("K61").Value = 6483
$ws.Range("L61").Value = 4426
$ws.Range("M61").Value = -6271
$ws.Range("N61").Value = -4850
$ws.Range("H88").Value = 1552.25
$ws.Range("J88").Value = 1429.4
$ws.Range("L88").Value = 1429.4
$ws.Range("N88").Value = -2241.4
$ws.Range("H91").Value = 1552.25
$ws.Range("J91").Value = 1429.4
$ws.Range("L91").Value = 1429.4
$ws.Range("N91").Value = -4237.4
$ws.Range("H102").Value = 864.6667
$ws.Range("I102").Value = 797.8182
$ws.Range("K102").Value = 797.8182
$ws.Range("M102").Value = 824.1818
$ws.Range("H110").Value = 948.6667
$ws.Range("I110").Value = 948.6667
$ws.Range("K110").Value = 948.6667
$ws.Range("M110").Value = 1096.3333
$ws.Range("H122").Value = 2755
$ws.Range("I122").Value = 2515.1667
$ws.Range("K122").Value = 7545.500100000001
$ws.Range("M122").Value = -5095.500100000001
$ws.Range("H136").Value = 4866.7856
$ws.Range("I136").Value = 6483
$ws.Range("J136").Value = 4426
$ws.Range("K136").Value = 19449
$ws.Range("L136").Value = 13278
$ws.Range("M136").Value = -16899
$ws.Range("N136").Value = -18378
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5899
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 5899
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 5899
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -9393
$ws.Range("H134").Value = 83337830
$ws.Range("I134").Value = 10000
$ws.Range("K134").Value = 30000
$ws.Range("M134").Value = -27465
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1965.7142
$ws.Range("I31").Value = 1785.0714
$ws.Range("K31").Value = 1785.0714
$ws.Range("M31").Value = -1490.0714
$ws.Range("H34").Value = 1965.7142
$ws.Range("I34").Value = 1785.0714
$ws.Range("K34").Value = 1785.0714
$ws.Range("M34").Value = -1583.0714
$ws.Range("H42").Value = 1885.3334
$ws.Range("I42").Value = 1885.3334
$ws.Range("K42").Value = 1885.3334
$ws.Range("M42").Value = -1292.3334
$ws.Range("H62").Value = 5476.0527
$ws.Range("I62").Value = 3682.1428
$ws.Range("J62").Value = 10499
$ws.Range("K62").Value = 3682.1428
$ws.Range("L62").Value = 10499
$ws.Range("M62").Value = -3058.1428
$ws.Range("N62").Value = -11747
$ws.Range("H65").Value = 5476.0527
$ws.Range("I65").Value = 3682.1428
$ws.Range("J65").Value = 10499
$ws.Range("K65").Value = 18410.714
$ws.Range("L65").Value = 52495
$ws.Range("M65").Value = -15290.714
$ws.Range("N65").Value = -58735
$ws.Range("H105").Value = 2599.6
$ws.Range("I105").Value = 1999.6666
$ws.Range("K105").Value = 1999.6666
$ws.Range("M105").Value = -252.6666
$ws.Range("H132").Value = 7929.5713
$ws.Range("I132").Value = 8284.362999999999
$ws.Range("J132").Value = 6628.6665
$ws.Range("K132").Value = 24853.089
$ws.Range("L132").Value = 19885.9995
$ws.Range("M132").Value = -22323.089
$ws.Range("N132").Value = -24945.9995
$ws.Range("H134").Value = 7146194
$ws.Range("I134").Value = 3318.25
$ws.Range("J134").Value = 50003450
$ws.Range("K134").Value = 9954.75
$ws.Range("L134").Value = 150010350
$ws.Range("M134").Value = -7419.75
$ws.Range("N134").Value = -150015420
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 687.7778
$ws.Range("I113").Value = 423
$ws.Range("K113").Value = 1269
$ws.Range("M113").Value = 901
$ws.Range("H122").Value = 1712.6923
$ws.Range("I122").Value = 923.5
$ws.Range("K122").Value = 8311.5
$ws.Range("M122").Value = -5861.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 496.66666
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 1873.6
$ws.Range("I132").Value = 2366
$ws.Range("J132").Value = 1135
$ws.Range("K132").Value = 7098
$ws.Range("L132").Value = 3405
$ws.Range("M132").Value = -4568
$ws.Range("N132").Value = -8465
$ws.Range("H136").Value = 38817.75
$ws.Range("I136").Value = 20999
$ws.Range("J136").Value = 44757.332
$ws.Range("K136").Value = 62997
$ws.Range("L136").Value = 134271.996
$ws.Range("M136").Value = -60447
$ws.Range("N136").Value = -139371.996
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2248.5
$ws.Range("I7").Value = 2248.5
$ws.Range("K7").Value = 2248.5
$ws.Range("M7").Value = -2136.5
$ws.Range("H16").Value = 733.375
$ws.Range("I16").Value = 733.375
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 733.375
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -563.375
$ws.Range("N16").ClearContents()
$ws.Range("H93").Value = 987.1667
$ws.Range("I93").Value = 847.2222
$ws.Range("K93").Value = 847.2222
$ws.Range("M93").Value = 400.7778
$ws.Range("H126").Value = 2248.5
$ws.Range("I126").Value = 2248.5
$ws.Range("K126").Value = 6745.5
$ws.Range("M126").Value = -4275.5
$ws.Range("H136").Value = 200004770
$ws.Range("I136").Value = 6934
$ws.Range("J136").Value = 333336670
$ws.Range("K136").Value = 20802
$ws.Range("L136").Value = 1000010010
$ws.Range("M136").Value = -18252
$ws.Range("N136").Value = -1000015110
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3199.3333
$ws.Range("I17").Value = 3199.3333
$ws.Range("K17").Value = 3199.3333
$ws.Range("M17").Value = -3027.3333
$ws.Range("H41").Value = 92335
$ws.Range("J41").Value = 98997.336
$ws.Range("L41").Value = 98997.336
$ws.Range("N41").Value = -99777.336
$ws.Range("H126").Value = 2609.4285
$ws.Range("I126").Value = 1759
$ws.Range("K126").Value = 5277
$ws.Range("M126").Value = -2807
